# Natmi following Dr Hou advice
# Update the LR-pairs data: new sending-cluster "M2" category inserted,
# a new sending cluster row-block for "sCs" appended, and the underlying
# expression/specificity values recomputed across all existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Reln", "Itgb1", "ECs", 3.0, 1.0, 0.2047723333333333, 0.614317, 0.04628602070509372, 0.04628602070509372, 3.0, 1.0, 153.5290173333333, 460.587052, 0.3172206968818489, 0.317220696881849, 31.43849511372044, 282.946456023484, 0.01468288374395752, 0.01468288374395752),
    @("ECs", "Reln", "Itgb1", "FAPs", 3.0, 1.0, 0.2047723333333333, 0.614317, 0.04628602070509372, 0.04628602070509372, 3.0, 1.0, 168.7997026666667, 506.3991080000001, 0.3487728915577651, 0.3487728915577651, 34.56550898102623, 311.0895808292361, 0.01614330928001812, 0.01614330928001812),
    @("ECs", "Reln", "Itgb1", "M2", 3.0, 1.0, 0.2047723333333333, 0.614317, 0.04628602070509372, 0.04628602070509372, 3.0, 1.0, 68.09032333333333, 204.27097, 0.1406878008722904, 0.1406878008722904, 13.94301438638778, 125.48712947749, 0.006511878464128935, 0.006511878464128935),
    @("ECs", "Reln", "Itgb1", "sCs", 3.0, 1.0, 0.2047723333333333, 0.614317, 0.04628602070509372, 0.04628602070509372, 3.0, 1.0, 93.562673, 280.688019, 0.1933186106880956, 0.1933186106880956, 19.15904686311367, 172.431421768023, 0.008947949216989146, 0.008947949216989146),
    @("FAPs", "Reln", "Itgb1", "ECs", 3.0, 1.0, 1.029733, 3.089199, 0.2327572391390028, 0.2327572391390028, 3.0, 1.0, 153.5290173333333, 460.587052, 0.3172206968818489, 0.317220696881849, 158.0938956057053, 1422.845060451348, 0.07383541360396964, 0.07383541360396965),
    @("FAPs", "Reln", "Itgb1", "FAPs", 3.0, 1.0, 1.029733, 3.089199, 0.2327572391390028, 0.2327572391390028, 3.0, 1.0, 168.7997026666667, 506.3991080000001, 0.3487728915577651, 0.3487728915577651, 173.8186242260547, 1564.367618034492, 0.08117941532551223, 0.08117941532551223),
    @("FAPs", "Reln", "Itgb1", "M2", 3.0, 1.0, 1.029733, 3.089199, 0.2327572391390028, 0.2327572391390028, 3.0, 1.0, 68.09032333333333, 204.27097, 0.1406878008722904, 0.1406878008722904, 70.11485291700333, 631.03367625303, 0.0327461041115721, 0.03274610411157211),
    @("FAPs", "Reln", "Itgb1", "sCs", 3.0, 1.0, 1.029733, 3.089199, 0.2327572391390028, 0.2327572391390028, 3.0, 1.0, 93.562673, 280.688019, 0.1933186106880956, 0.1933186106880956, 96.34457195630901, 867.101147606781, 0.04499630609794886, 0.04499630609794886),
    @("M2", "Reln", "Itgb1", "ECs", 1.0, 0.3333333333333333, 0.004225666666666667, 0.012677, 0.0009551548866114287, 0.0009551548866114289, 3.0, 1.0, 153.5290173333333, 460.587052, 0.3172206968818489, 0.317220696881849, 0.6487624509115555, 5.838862058204, 0.0003029948987609808, 0.0003029948987609809),
    @("M2", "Reln", "Itgb1", "FAPs", 1.0, 0.3333333333333333, 0.004225666666666667, 0.012677, 0.0009551548866114287, 0.0009551548866114289, 3.0, 1.0, 168.7997026666667, 506.3991080000001, 0.3487728915577651, 0.3487728915577651, 0.7132912769017778, 6.419621492116002, 0.0003331321316889973, 0.0003331321316889973),
    @("M2", "Reln", "Itgb1", "M2", 1.0, 0.3333333333333333, 0.004225666666666667, 0.012677, 0.0009551548866114287, 0.0009551548866114289, 3.0, 1.0, 68.09032333333333, 204.27097, 0.1406878008722904, 0.1406878008722904, 0.2877270096322222, 2.58954308669, 0.0001343786404897838, 0.0001343786404897838),
    @("M2", "Reln", "Itgb1", "sCs", 1.0, 0.3333333333333333, 0.004225666666666667, 0.012677, 0.0009551548866114287, 0.0009551548866114289, 3.0, 1.0, 93.562673, 280.688019, 0.1933186106880956, 0.1933186106880956, 0.3953646685403334, 3.558282016863, 0.0001846492156716669, 0.0001846492156716669),
    @("sCs", "Reln", "Itgb1", "ECs", 3.0, 1.0, 3.185333333333334, 9.556000000000001, 0.720001585269292, 0.7200015852692919, 3.0, 1.0, 153.5290173333333, 460.587052, 0.3172206968818489, 0.317220696881849, 489.0410965457778, 4401.369868912, 0.2283994046351608, 0.2283994046351608),
    @("sCs", "Reln", "Itgb1", "FAPs", 3.0, 1.0, 3.185333333333334, 9.556000000000001, 0.720001585269292, 0.7200015852692919, 3.0, 1.0, 168.7997026666667, 506.3991080000001, 0.3487728915577651, 0.3487728915577651, 537.683319560889, 4839.149876048001, 0.2511170348205458, 0.2511170348205457),
    @("sCs", "Reln", "Itgb1", "M2", 3.0, 1.0, 3.185333333333334, 9.556000000000001, 0.720001585269292, 0.7200015852692919, 3.0, 1.0, 68.09032333333333, 204.27097, 0.1406878008722904, 0.1406878008722904, 216.8903765911111, 1952.01338932, 0.1012954396560996, 0.1012954396560996),
    @("sCs", "Reln", "Itgb1", "sCs", 3.0, 1.0, 3.185333333333334, 9.556000000000001, 0.720001585269292, 0.7200015852692919, 3.0, 1.0, 93.562673, 280.688019, 0.1933186106880956, 0.1933186106880956, 298.0283010626667, 2682.254709564, 0.139189706157486, 0.1391897061574859)
)


$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
